$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.381914138793945
$ws.Range("B1").Value = 2.658785343170166
$ws.Range("C1").Value = 5.866413593292236
$ws.Range("D1").Value = 2.266639232635498
$ws.Range("E1").Value = 1.210982084274292
